# Update column F ("dSF") values for a subset of rows.
# These values were repulled from source data and no longer equal column E ("dS0").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    7  = -1
    10 = -1
    14 = -2
    16 = 1
    22 = -2
    25 = -2
    26 = -4
    30 = -3
    32 = -5
    33 = 3
    34 = 6
    38 = 4
    40 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
